$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").Value = "63.418.96"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.672.76"
$ws.Range("E3").Value = "  +3.81%  "

$ws.Range("D5").Value = "'611.17"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("D6").Value = "'143.65"
$ws.Range("E6").Value = "  -0.62%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").Value = "2.669.59"
$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").Value = "'5.62"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D13").Value = "'0.363"
$ws.Range("E13").Value = "  +3.54%  "

$ws.Range("D14").Value = "'27.36"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "3.148.85"

$ws.Range("D16").Value = "63.245.99"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "2.677.37"
$ws.Range("E18").Value = "  +4.17%  "

$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  +3.22%  "

$ws.Range("D20").Value = "'341.77"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'4.42"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("D22").Value = "'6.90"
$ws.Range("E22").Value = "  +3.79%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'67.14"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("E25").Value = "  +1.87%  "

$ws.Range("E26").Value = "  -2.47%  "

$ws.Range("D27").Value = "'8.66"
$ws.Range("E27").Value = "  +4.95%  "

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").Value = "'544.28"
$ws.Range("E29").Value = "  +15.95%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("E32").Value = "  +5.73%  "

$ws.Range("E33").Value = "  +7.29%  "

$ws.Range("D34").Value = "'" + "0.0" + $sub3 + "0807"
$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D35").Value = "'172.98"
$ws.Range("E35").Value = "  -2.00%  "

$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "  +14.39%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +1.27%  "

$ws.Range("D39").Value = "'19.19"
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("D40").Value = "'1.87"
$ws.Range("E40").Value = "  +10.00%  "

$ws.Range("E41").Value = "  +12.37%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("D44").Value = "'22.28"
$ws.Range("E44").Value = "  +4.16%  "

$ws.Range("D45").Value = "'0.0575"
$ws.Range("E45").Value = "  +6.80%  "

$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("E47").Value = "  +1.47%  "

$ws.Range("D48").Value = "'0.0963"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'18.74"
$ws.Range("E49").Value = "  +3.13%  "

$ws.Range("E50").Value = "  +4.48%  "

$ws.Range("D51").Value = "'11.30"
$ws.Range("E51").Value = "  -0.76%  "
